$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2174
$ws1.Range("F4").Value = 48
$ws1.Range("F5").Value = 11475
$ws1.Range("F9").Value = 11423
$ws1.Range("F10").Value = 466
$ws1.Range("F12").Value = 76
$ws1.Range("F14").Value = 5680
$ws1.Range("F16").Value = 3487

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2174
$ws4.Range("F5").Value = 48
$ws4.Range("F7").Value = 11475
$ws4.Range("F11").Value = 11423
$ws4.Range("F12").Value = 466
$ws4.Range("F14").Value = 76
$ws4.Range("F17").Value = 5680
$ws4.Range("F19").Value = 3487
